$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 12.9
$ws.Range("B3").Value = 4.5999999999999996
$ws.Range("C3").Value = 11.7
$ws.Range("B4").Value = 0.65
